$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = ""
$ws.Cells.Item(6,6).Value = ""
$ws.Cells.Item(6,5).Value = "Started front-end work regarding the main page through following tutorials suggested by the team. Helped explain the models within the feasibility study and did some proof reading/editing. Finished the heuristics document for user experience. Getting used to the new tools, so far so good. "
$ws.Range("E6").Select()
